$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A198").Value = "TAO-USD"
$ws.Range("A199").Value = "IMX-USD"
$ws.Range("A200").Value = "GRT-USD"
